# "Generate Report for Archive"
# - The localization status for the (still in-flight) e2e files moved on from
#   "Ready for handoff" to "In Translation" (shows up on the Overview sheet's
#   per-locale status columns as well as each locale sheet's "Status" column
#   -- they all share the same underlying string).
# - Regenerating the report re-sized the now-narrower "Status"/locale-status
#   columns to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

# --- 1) Flip every "Ready for handoff" cell to "In Translation" -----------
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($cell.Value2 -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# --- 2) Re-fit the columns that held that text to their new width ---------
# ColumnWidth is expressed in characters and Excel snaps it to whole screen
# pixels, so we dial in the input that lands on the narrowest/closest pixel
# column to the regenerated report's target width.
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Columns.Item(5).ColumnWidth = 12.5   # zh-cn status column
$ws1.Columns.Item(6).ColumnWidth = 12.5   # de-de status column

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Columns.Item(3).ColumnWidth = 12.5   # Status column

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Columns.Item(3).ColumnWidth = 12.5   # Status column
